$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('CJ McCollum', 'PG,SG', 'New Orleans Pelicans'),
    @('Jalen Williams', 'SG,SF,PF,C', 'Oklahoma City Thunder'),
    @('Shai Gilgeous-Alexander', 'PG,SG', 'Oklahoma City Thunder'),
    @('Lauri Markkanen', 'SF,PF', 'Utah Jazz'),
    @('Jordan Clarkson', 'SG,SF', 'Utah Jazz'),
    @('Jordan Poole', 'PG,SG', 'Washington Wizards'),
    @('Keyonte George', 'PG,SG', 'Utah Jazz'),
    @('Jimmy Butler', 'SF,PF', 'Miami Heat'),
    @('Joel Embiid', 'C', 'Philadelphia 76ers'),
    @('Christian Braun', 'SG,SF', 'Denver Nuggets'),
    @('Tobias Harris', 'SF,PF', 'Detroit Pistons'),
    @('Zach LaVine', 'SG,SF', 'Chicago Bulls'),
    @('Andrew Nembhard', 'PG,SG', 'Indiana Pacers'),
    @('Guerschon Yabusele', 'PF,C', 'Philadelphia 76ers'),
    @('John Collins', 'PF,C', 'Utah Jazz'),
    @('RJ Barrett', 'SG,SF,PF', 'Toronto Raptors'),
    @('Kyrie Irving', 'PG,SG', 'Dallas Mavericks')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
